# Update cryptocurrency price/volume data per upstream diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.264.17"
$ws.Range("E2").Value = "  +2.55%  "
$ws.Range("D3").Value = "'2.970.11"
$ws.Range("E3").Value = "  +0.93%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'596.84"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").Value = "'149.97"
$ws.Range("E6").Value = "  +2.88%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "'2.970.09"
$ws.Range("E8").Value = "  +0.95%  "
$ws.Range("E9").Value = "  +1.04%  "
$ws.Range("D10").Value = "'7.35"
$ws.Range("E10").Value = "  +5.39%  "
$ws.Range("E11").Value = "  +8.66%  "
$ws.Range("D12").Value = "'0.448"
$ws.Range("E12").Value = "  +1.38%  "
$ws.Range("D13").Value = "'0.0000245"
$ws.Range("E13").Value = "  +8.18%  "
$ws.Range("D14").Value = "'33.38"
$ws.Range("E14").Value = "  -1.27%  "
$ws.Range("E15").Value = "  -0.59%  "
$ws.Range("D16").Value = "'3.464.46"
$ws.Range("E16").Value = "  +1.15%  "
$ws.Range("D17").Value = "'63.164.72"
$ws.Range("E17").Value = "  +2.43%  "
$ws.Range("D18").Value = "'6.81"
$ws.Range("E18").Value = "  +0.99%  "
$ws.Range("D19").Value = "'2.969.68"
$ws.Range("E19").Value = "  +0.86%  "
$ws.Range("D20").Value = "'446.63"
$ws.Range("E20").Value = "  +2.80%  "
$ws.Range("D21").Value = "'13.63"
$ws.Range("E21").Value = "  +0.67%  "
$ws.Range("D22").Value = "'0.675"
$ws.Range("E22").Value = "  -0.91%  "
$ws.Range("D23").Value = "'7.17"
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").Value = "'11.44"
$ws.Range("E24").Value = "  +4.04%  "
$ws.Range("D25").Value = "'81.99"
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("D26").Value = "'2.21"
$ws.Range("E26").Value = "  -0.71%  "
$ws.Range("D27").Value = "'11.97"
$ws.Range("E27").Value = "  +0.71%  "
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("D29").Value = "'0.0000108"
$ws.Range("E29").Value = "  +22.43%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "'2.25"
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "'7.32"
$ws.Range("E31").Value = "  +4.55%  "
$ws.Range("E32").Value = "  +0.93%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "'26.90"
$ws.Range("E33").Value = "  +0.46%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.111"
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("E35").Value = "  -0.18%  "
$ws.Range("B36").Value = "dogwifhat"
$ws.Range("C36").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D36").Value = "'3.33"
$ws.Range("E36").Value = "  +9.87%  "
$ws.Range("B37").Value = "Mantle"
$ws.Range("C37").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  -1.58%  "
$ws.Range("D38").Value = "'5.70"
$ws.Range("E38").Value = "  +0.70%  "
$ws.Range("E39").Value = "  +2.95%  "
$ws.Range("D40").Value = "'49.81"
$ws.Range("E40").Value = "  -0.41%  "
$ws.Range("D41").Value = "'8.62"
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("E42").Value = "  -4.44%  "
$ws.Range("D43").Value = "'0.287"
$ws.Range("E43").Value = "  +0.88%  "
$ws.Range("D44").Value = "'41.40"
$ws.Range("E44").Value = "  -2.84%  "
$ws.Range("D45").Value = "'2.719.51"
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D46").Value = "'371.72"
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "'0.0344"
$ws.Range("E47").Value = "  -1.73%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").Value = "'135.17"
$ws.Range("E48").Value = "  +1.06%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").Value = "'23.40"
$ws.Range("E50").Value = "  -1.70%  "
$ws.Range("D51").Value = "'0.106"
$ws.Range("E51").Value = "  -0.14%  "
